$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Value = "2025-04-29 00:45:21"
$ws.Range("B37").Value = 128
